# SettingSharedFormula example workbook — refresh the Aspose.Cells
# evaluation-version watermark that lives on the "Evaluation Warning"
# sheet (this is what changes when the sample is regenerated with a
# newer library build: "...Copyright 2003 - 2014..." becomes
# "...Copyright 2003 - 2016...").
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Evaluation Warning")
$ws.Range("A5").Value = "Evaluation Only. Created with Aspose.Cells for .NET.Copyright 2003 - 2016 Aspose Pty Ltd."
